$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 700.4545000000001
$ws.Range("I39").Value = 533.5
$ws.Range("K39").Value = 1600.5
$ws.Range("M39").Value = -1304.5
$ws.Range("H64").Value = 5142.8423
$ws.Range("I64").Value = 4153.4614
$ws.Range("J64").Value = 7286.5
$ws.Range("K64").Value = 4153.4614
$ws.Range("L64").Value = 7286.5
$ws.Range("M64").Value = -3905.4614
$ws.Range("N64").Value = -7782.5
$ws.Range("H67").Value = 5142.8423
$ws.Range("I67").Value = 4153.4614
$ws.Range("J67").Value = 7286.5
$ws.Range("K67").Value = 4153.4614
$ws.Range("L67").Value = 7286.5
$ws.Range("M67").Value = -3295.4614
$ws.Range("N67").Value = -9002.5
$ws.Range("H70").Value = 10010.04
$ws.Range("J70").Value = 13358.866
$ws.Range("L70").Value = 40076.598
$ws.Range("N70").Value = -40616.598
$ws.Range("H73").Value = 10010.04
$ws.Range("J73").Value = 13358.866
$ws.Range("L73").Value = 40076.598
$ws.Range("N73").Value = -41948.598
$ws.Range("H86").Value = 115449.78
$ws.Range("I86").Value = 115449.78
$ws.Range("K86").Value = 115449.78
$ws.Range("M86").Value = -114326.78
$ws.Range("H89").Value = 115449.78
$ws.Range("I89").Value = 115449.78
$ws.Range("K89").Value = 577248.9
$ws.Range("M89").Value = -571632.9
$ws.Range("H100").Value = 1440.7894
$ws.Range("I100").Value = 1109.7693
$ws.Range("J100").Value = 2158
$ws.Range("K100").Value = 1109.7693
$ws.Range("L100").Value = 2158
$ws.Range("M100").Value = -568.7692999999999
$ws.Range("N100").Value = -3240
$ws.Range("H116").Value = 8565.559999999999
$ws.Range("I116").Value = 9958.058999999999
$ws.Range("J116").Value = 5606.5
$ws.Range("K116").Value = 9958.058999999999
$ws.Range("L116").Value = 5606.5
$ws.Range("M116").Value = -6516.058999999999
$ws.Range("N116").Value = -12490.5
$ws.Range("H135").Value = 2619.6843
$ws.Range("J135").Value = 6896.75
$ws.Range("L135").Value = 62070.75
$ws.Range("N135").Value = -67140.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 4852
$ws.Range("I26").Value = 4450
$ws.Range("K26").Value = 4450
$ws.Range("M26").Value = -4120
$ws.Range("H32").Value = 7582.3657
$ws.Range("I32").Value = 7203.148
$ws.Range("K32").Value = 7203.148
$ws.Range("M32").Value = -6916.148
$ws.Range("H61").Value = 7095.25
$ws.Range("I61").Value = 2000.7693
$ws.Range("K61").Value = 2000.7693
$ws.Range("M61").Value = -1788.7693
$ws.Range("H132").Value = 3272.561
$ws.Range("I132").Value = 3266.925
$ws.Range("K132").Value = 9800.775000000001
$ws.Range("M132").Value = -7270.775000000001
$ws.Range("H136").Value = 7095.25
$ws.Range("I136").Value = 2000.7693
$ws.Range("K136").Value = 6002.3079
$ws.Range("M136").Value = -3452.3079
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 48499
$ws.Range("J40").Value = 48499
$ws.Range("L40").Value = 48499
$ws.Range("N40").Value = -49029
$ws.Range("H107").Value = 2145.25
$ws.Range("I107").Value = 1100.6666
$ws.Range("K107").Value = 1100.6666
$ws.Range("M107").Value = 819.3334
$ws.Range("H134").Value = 1963.1
$ws.Range("I134").Value = 1230.9656
$ws.Range("J134").Value = 2974.1428
$ws.Range("K134").Value = 3692.8968
$ws.Range("L134").Value = 8922.428400000001
$ws.Range("M134").Value = -1157.8968
$ws.Range("N134").Value = -13992.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 30000
$ws.Range("K6").Value = 30000
$ws.Range("M6").Value = -29887
$ws.Range("H86").Value = 4800.1113
$ws.Range("I86").Value = 4016.8
$ws.Range("K86").Value = 4016.8
$ws.Range("M86").Value = -2893.8
$ws.Range("H89").Value = 4800.1113
$ws.Range("I89").Value = 4016.8
$ws.Range("K89").Value = 20084
$ws.Range("M89").Value = -14468
$ws.Range("H99").Value = 436855.88
$ws.Range("I99").Value = 557566.7
$ws.Range("K99").Value = 557566.7
$ws.Range("M99").Value = -556068.7
$ws.Range("H126").Value = 436855.88
$ws.Range("I126").Value = 557566.7
$ws.Range("K126").Value = 1672700.1
$ws.Range("M126").Value = -1670230.1
$ws.Range("H132").Value = 4156.32
$ws.Range("I132").Value = 2278.6086
$ws.Range("J132").Value = 25750
$ws.Range("K132").Value = 6835.825800000001
$ws.Range("L132").Value = 77250
$ws.Range("M132").Value = -4305.825800000001
$ws.Range("N132").Value = -82310
$ws.Range("H134").Value = 2866.7585
$ws.Range("I134").Value = 2804.818
$ws.Range("J134").Value = 4002.3333
$ws.Range("K134").Value = 8414.454000000002
$ws.Range("L134").Value = 12006.9999
$ws.Range("M134").Value = -5879.454000000002
$ws.Range("N134").Value = -17076.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 314754.16
$ws.Range("I4").Value = 486178.38
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 1458535.14
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = -1458423.14
$ws.Range("N4").Value = -30224
$ws.Range("H5").Value = 1140.3
$ws.Range("I5").Value = 1140.3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3420.9
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3308.9
$ws.Range("N5").ClearContents()
$ws.Range("H131").Value = 11629404
$ws.Range("J131").Value = 1554.375
$ws.Range("L131").Value = 4663.125
$ws.Range("N131").Value = -14743.125
$ws.Range("H134").Value = 5968.3335
$ws.Range("I134").Value = 5968.3335
$ws.Range("K134").Value = 17905.0005
$ws.Range("M134").Value = -12835.0005
$ws.Range("H135").Value = 1140.3
$ws.Range("I135").Value = 1140.3
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10262.7
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7727.699999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 32350
$ws.Range("J104").Value = 32350
$ws.Range("L104").Value = 32350
$ws.Range("N104").Value = -39338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2344.9285
$ws.Range("I68").Value = 2302.3845
$ws.Range("J68").Value = 2898
$ws.Range("K68").Value = 2302.3845
$ws.Range("L68").Value = 2898
$ws.Range("M68").Value = -1553.3845
$ws.Range("N68").Value = -4396
$ws.Range("H71").Value = 2344.9285
$ws.Range("I71").Value = 2302.3845
$ws.Range("J71").Value = 2898
$ws.Range("K71").Value = 11511.9225
$ws.Range("L71").Value = 14490
$ws.Range("M71").Value = -7767.922500000001
$ws.Range("N71").Value = -21978
$ws.Range("H82").Value = 2611.2222
$ws.Range("I82").Value = 2312.25
$ws.Range("K82").Value = 2312.25
$ws.Range("M82").Value = -1951.25
$ws.Range("H85").Value = 2611.2222
$ws.Range("I85").Value = 2312.25
$ws.Range("K85").Value = 2312.25
$ws.Range("M85").Value = -1064.25
$ws.Range("H136").Value = 7144.385
$ws.Range("I136").Value = 3573.0833
$ws.Range("J136").Value = 50000
$ws.Range("K136").Value = 10719.2499
$ws.Range("L136").Value = 150000
$ws.Range("M136").Value = -8169.249899999999
$ws.Range("N136").Value = -155100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11509.375
$ws.Range("J62").Value = 11867.857
$ws.Range("L62").Value = 11867.857
$ws.Range("N62").Value = -13115.857
$ws.Range("H65").Value = 11509.375
$ws.Range("J65").Value = 11867.857
$ws.Range("L65").Value = 59339.285
$ws.Range("N65").Value = -65579.285
$ws.Range("H132").Value = 2285.75
$ws.Range("I132").Value = 1214.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3643.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1113.5
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 235695.62
$ws.Range("I136").Value = 288465.56
$ws.Range("J136").Value = 4827.125
$ws.Range("K136").Value = 865396.6799999999
$ws.Range("L136").Value = 14481.375
$ws.Range("M136").Value = -862846.6799999999
$ws.Range("N136").Value = -19581.375

